$wb = $excel.ActiveWorkbook

# --- Sheet: ALERTS (rows 17-17) ---
$ws = $wb.Worksheets.Item('ALERTS')
# Row 17
$ws.Cells.Item(17, 1).NumberFormat = '@'
$ws.Cells.Item(17, 1).Value = '2026-01-30'
$ws.Cells.Item(17, 1).ClearFormats()
$ws.Cells.Item(17, 2).Value = '17:52:12'
$ws.Cells.Item(17, 3).Value = '17:00'
$ws.Cells.Item(17, 4).Value = 'Bathroom'
$ws.Cells.Item(17, 5).Value = 'MINIMAL'
$ws.Cells.Item(17, 6).Value = 'MINIMAL ALERT: Bathroom occupied, no motion > 20s.'

# --- Sheet: PIR (rows 357-369) ---
$ws = $wb.Worksheets.Item('PIR')
# Row 357
$ws.Cells.Item(357, 1).NumberFormat = '@'
$ws.Cells.Item(357, 1).Value = '2026-01-30'
$ws.Cells.Item(357, 1).ClearFormats()
$ws.Cells.Item(357, 2).Value = '17:51:50'
$ws.Cells.Item(357, 3).Value = '17:00'
$ws.Cells.Item(357, 4).Value = 'Bathroom'
$ws.Cells.Item(357, 5).Value = 'No Motion'
$ws.Cells.Item(357, 6).Value = 'Inactive'
# Row 358
$ws.Cells.Item(358, 1).NumberFormat = '@'
$ws.Cells.Item(358, 1).Value = '2026-01-30'
$ws.Cells.Item(358, 1).ClearFormats()
$ws.Cells.Item(358, 2).Value = '17:51:52'
$ws.Cells.Item(358, 3).Value = '17:00'
$ws.Cells.Item(358, 4).Value = 'Bathroom'
$ws.Cells.Item(358, 5).Value = 'No Motion'
$ws.Cells.Item(358, 6).Value = 'Inactive'
# Row 359
$ws.Cells.Item(359, 1).NumberFormat = '@'
$ws.Cells.Item(359, 1).Value = '2026-01-30'
$ws.Cells.Item(359, 1).ClearFormats()
$ws.Cells.Item(359, 2).Value = '17:51:56'
$ws.Cells.Item(359, 3).Value = '17:00'
$ws.Cells.Item(359, 4).Value = 'Bathroom'
$ws.Cells.Item(359, 5).Value = 'No Motion'
$ws.Cells.Item(359, 6).Value = 'Inactive'
# Row 360
$ws.Cells.Item(360, 1).NumberFormat = '@'
$ws.Cells.Item(360, 1).Value = '2026-01-30'
$ws.Cells.Item(360, 1).ClearFormats()
$ws.Cells.Item(360, 2).Value = '17:52:01'
$ws.Cells.Item(360, 3).Value = '17:00'
$ws.Cells.Item(360, 4).Value = 'Bathroom'
$ws.Cells.Item(360, 5).Value = 'No Motion'
$ws.Cells.Item(360, 6).Value = 'Inactive'
# Row 361
$ws.Cells.Item(361, 1).NumberFormat = '@'
$ws.Cells.Item(361, 1).Value = '2026-01-30'
$ws.Cells.Item(361, 1).ClearFormats()
$ws.Cells.Item(361, 2).Value = '17:52:06'
$ws.Cells.Item(361, 3).Value = '17:00'
$ws.Cells.Item(361, 4).Value = 'Bathroom'
$ws.Cells.Item(361, 5).Value = 'No Motion'
$ws.Cells.Item(361, 6).Value = 'Inactive'
# Row 362
$ws.Cells.Item(362, 1).NumberFormat = '@'
$ws.Cells.Item(362, 1).Value = '2026-01-30'
$ws.Cells.Item(362, 1).ClearFormats()
$ws.Cells.Item(362, 2).Value = '17:52:11'
$ws.Cells.Item(362, 3).Value = '17:00'
$ws.Cells.Item(362, 4).Value = 'Bathroom'
$ws.Cells.Item(362, 5).Value = 'No Motion'
$ws.Cells.Item(362, 6).Value = 'Inactive'
# Row 363
$ws.Cells.Item(363, 1).NumberFormat = '@'
$ws.Cells.Item(363, 1).Value = '2026-01-30'
$ws.Cells.Item(363, 1).ClearFormats()
$ws.Cells.Item(363, 2).Value = '17:52:16'
$ws.Cells.Item(363, 3).Value = '17:00'
$ws.Cells.Item(363, 4).Value = 'Bathroom'
$ws.Cells.Item(363, 5).Value = 'No Motion'
$ws.Cells.Item(363, 6).Value = 'Inactive'
# Row 364
$ws.Cells.Item(364, 1).NumberFormat = '@'
$ws.Cells.Item(364, 1).Value = '2026-01-30'
$ws.Cells.Item(364, 1).ClearFormats()
$ws.Cells.Item(364, 2).Value = '17:52:21'
$ws.Cells.Item(364, 3).Value = '17:00'
$ws.Cells.Item(364, 4).Value = 'Bathroom'
$ws.Cells.Item(364, 5).Value = 'No Motion'
$ws.Cells.Item(364, 6).Value = 'Inactive'
# Row 365
$ws.Cells.Item(365, 1).NumberFormat = '@'
$ws.Cells.Item(365, 1).Value = '2026-01-30'
$ws.Cells.Item(365, 1).ClearFormats()
$ws.Cells.Item(365, 2).Value = '17:52:26'
$ws.Cells.Item(365, 3).Value = '17:00'
$ws.Cells.Item(365, 4).Value = 'Bathroom'
$ws.Cells.Item(365, 5).Value = 'No Motion'
$ws.Cells.Item(365, 6).Value = 'Inactive'
# Row 366
$ws.Cells.Item(366, 1).NumberFormat = '@'
$ws.Cells.Item(366, 1).Value = '2026-01-30'
$ws.Cells.Item(366, 1).ClearFormats()
$ws.Cells.Item(366, 2).Value = '17:52:31'
$ws.Cells.Item(366, 3).Value = '17:00'
$ws.Cells.Item(366, 4).Value = 'Bathroom'
$ws.Cells.Item(366, 5).Value = 'No Motion'
$ws.Cells.Item(366, 6).Value = 'Inactive'
# Row 367
$ws.Cells.Item(367, 1).NumberFormat = '@'
$ws.Cells.Item(367, 1).Value = '2026-01-30'
$ws.Cells.Item(367, 1).ClearFormats()
$ws.Cells.Item(367, 2).Value = '17:52:37'
$ws.Cells.Item(367, 3).Value = '17:00'
$ws.Cells.Item(367, 4).Value = 'Bathroom'
$ws.Cells.Item(367, 5).Value = 'No Motion'
$ws.Cells.Item(367, 6).Value = 'Inactive'
# Row 368
$ws.Cells.Item(368, 1).NumberFormat = '@'
$ws.Cells.Item(368, 1).Value = '2026-01-30'
$ws.Cells.Item(368, 1).ClearFormats()
$ws.Cells.Item(368, 2).Value = '17:52:41'
$ws.Cells.Item(368, 3).Value = '17:00'
$ws.Cells.Item(368, 4).Value = 'Bathroom'
$ws.Cells.Item(368, 5).Value = 'No Motion'
$ws.Cells.Item(368, 6).Value = 'Inactive'
# Row 369
$ws.Cells.Item(369, 1).NumberFormat = '@'
$ws.Cells.Item(369, 1).Value = '2026-01-30'
$ws.Cells.Item(369, 1).ClearFormats()
$ws.Cells.Item(369, 2).Value = '17:52:46'
$ws.Cells.Item(369, 3).Value = '17:00'
$ws.Cells.Item(369, 4).Value = 'Bathroom'
$ws.Cells.Item(369, 5).Value = 'No Motion'
$ws.Cells.Item(369, 6).Value = 'Inactive'

# --- Sheet: Humidity (rows 244-254) ---
$ws = $wb.Worksheets.Item('Humidity')
# Row 244
$ws.Cells.Item(244, 1).NumberFormat = '@'
$ws.Cells.Item(244, 1).Value = '2026-01-30'
$ws.Cells.Item(244, 1).ClearFormats()
$ws.Cells.Item(244, 2).Value = '17:51:49'
$ws.Cells.Item(244, 3).Value = '17:00'
$ws.Cells.Item(244, 4).Value = 'Bathroom'
$ws.Cells.Item(244, 5).NumberFormat = '@'
$ws.Cells.Item(244, 5).Value = '85.6%'
$ws.Cells.Item(244, 5).ClearFormats()
$ws.Cells.Item(244, 6).Value = 'Active'
# Row 245
$ws.Cells.Item(245, 1).NumberFormat = '@'
$ws.Cells.Item(245, 1).Value = '2026-01-30'
$ws.Cells.Item(245, 1).ClearFormats()
$ws.Cells.Item(245, 2).Value = '17:51:51'
$ws.Cells.Item(245, 3).Value = '17:00'
$ws.Cells.Item(245, 4).Value = 'Bathroom'
$ws.Cells.Item(245, 5).NumberFormat = '@'
$ws.Cells.Item(245, 5).Value = '86.1%'
$ws.Cells.Item(245, 5).ClearFormats()
$ws.Cells.Item(245, 6).Value = 'Active'
# Row 246
$ws.Cells.Item(246, 1).NumberFormat = '@'
$ws.Cells.Item(246, 1).Value = '2026-01-30'
$ws.Cells.Item(246, 1).ClearFormats()
$ws.Cells.Item(246, 2).Value = '17:51:57'
$ws.Cells.Item(246, 3).Value = '17:00'
$ws.Cells.Item(246, 4).Value = 'Bathroom'
$ws.Cells.Item(246, 5).NumberFormat = '@'
$ws.Cells.Item(246, 5).Value = '87.0%'
$ws.Cells.Item(246, 5).ClearFormats()
$ws.Cells.Item(246, 6).Value = 'Active'
# Row 247
$ws.Cells.Item(247, 1).NumberFormat = '@'
$ws.Cells.Item(247, 1).Value = '2026-01-30'
$ws.Cells.Item(247, 1).ClearFormats()
$ws.Cells.Item(247, 2).Value = '17:52:02'
$ws.Cells.Item(247, 3).Value = '17:00'
$ws.Cells.Item(247, 4).Value = 'Bathroom'
$ws.Cells.Item(247, 5).NumberFormat = '@'
$ws.Cells.Item(247, 5).Value = '87.0%'
$ws.Cells.Item(247, 5).ClearFormats()
$ws.Cells.Item(247, 6).Value = 'Active'
# Row 248
$ws.Cells.Item(248, 1).NumberFormat = '@'
$ws.Cells.Item(248, 1).Value = '2026-01-30'
$ws.Cells.Item(248, 1).ClearFormats()
$ws.Cells.Item(248, 2).Value = '17:52:07'
$ws.Cells.Item(248, 3).Value = '17:00'
$ws.Cells.Item(248, 4).Value = 'Bathroom'
$ws.Cells.Item(248, 5).NumberFormat = '@'
$ws.Cells.Item(248, 5).Value = '86.0%'
$ws.Cells.Item(248, 5).ClearFormats()
$ws.Cells.Item(248, 6).Value = 'Active'
# Row 249
$ws.Cells.Item(249, 1).NumberFormat = '@'
$ws.Cells.Item(249, 1).Value = '2026-01-30'
$ws.Cells.Item(249, 1).ClearFormats()
$ws.Cells.Item(249, 2).Value = '17:52:17'
$ws.Cells.Item(249, 3).Value = '17:00'
$ws.Cells.Item(249, 4).Value = 'Bathroom'
$ws.Cells.Item(249, 5).NumberFormat = '@'
$ws.Cells.Item(249, 5).Value = '86.9%'
$ws.Cells.Item(249, 5).ClearFormats()
$ws.Cells.Item(249, 6).Value = 'Active'
# Row 250
$ws.Cells.Item(250, 1).NumberFormat = '@'
$ws.Cells.Item(250, 1).Value = '2026-01-30'
$ws.Cells.Item(250, 1).ClearFormats()
$ws.Cells.Item(250, 2).Value = '17:52:22'
$ws.Cells.Item(250, 3).Value = '17:00'
$ws.Cells.Item(250, 4).Value = 'Bathroom'
$ws.Cells.Item(250, 5).NumberFormat = '@'
$ws.Cells.Item(250, 5).Value = '87.0%'
$ws.Cells.Item(250, 5).ClearFormats()
$ws.Cells.Item(250, 6).Value = 'Active'
# Row 251
$ws.Cells.Item(251, 1).NumberFormat = '@'
$ws.Cells.Item(251, 1).Value = '2026-01-30'
$ws.Cells.Item(251, 1).ClearFormats()
$ws.Cells.Item(251, 2).Value = '17:52:27'
$ws.Cells.Item(251, 3).Value = '17:00'
$ws.Cells.Item(251, 4).Value = 'Bathroom'
$ws.Cells.Item(251, 5).NumberFormat = '@'
$ws.Cells.Item(251, 5).Value = '86.0%'
$ws.Cells.Item(251, 5).ClearFormats()
$ws.Cells.Item(251, 6).Value = 'Active'
# Row 252
$ws.Cells.Item(252, 1).NumberFormat = '@'
$ws.Cells.Item(252, 1).Value = '2026-01-30'
$ws.Cells.Item(252, 1).ClearFormats()
$ws.Cells.Item(252, 2).Value = '17:52:37'
$ws.Cells.Item(252, 3).Value = '17:00'
$ws.Cells.Item(252, 4).Value = 'Bathroom'
$ws.Cells.Item(252, 5).NumberFormat = '@'
$ws.Cells.Item(252, 5).Value = '86.0%'
$ws.Cells.Item(252, 5).ClearFormats()
$ws.Cells.Item(252, 6).Value = 'Active'
# Row 253
$ws.Cells.Item(253, 1).NumberFormat = '@'
$ws.Cells.Item(253, 1).Value = '2026-01-30'
$ws.Cells.Item(253, 1).ClearFormats()
$ws.Cells.Item(253, 2).Value = '17:52:42'
$ws.Cells.Item(253, 3).Value = '17:00'
$ws.Cells.Item(253, 4).Value = 'Bathroom'
$ws.Cells.Item(253, 5).NumberFormat = '@'
$ws.Cells.Item(253, 5).Value = '86.9%'
$ws.Cells.Item(253, 5).ClearFormats()
$ws.Cells.Item(253, 6).Value = 'Active'
# Row 254
$ws.Cells.Item(254, 1).NumberFormat = '@'
$ws.Cells.Item(254, 1).Value = '2026-01-30'
$ws.Cells.Item(254, 1).ClearFormats()
$ws.Cells.Item(254, 2).Value = '17:52:47'
$ws.Cells.Item(254, 3).Value = '17:00'
$ws.Cells.Item(254, 4).Value = 'Bathroom'
$ws.Cells.Item(254, 5).NumberFormat = '@'
$ws.Cells.Item(254, 5).Value = '85.9%'
$ws.Cells.Item(254, 5).ClearFormats()
$ws.Cells.Item(254, 6).Value = 'Active'

# --- Sheet: Proximity (rows 57-59) ---
$ws = $wb.Worksheets.Item('Proximity')
# Row 57
$ws.Cells.Item(57, 1).NumberFormat = '@'
$ws.Cells.Item(57, 1).Value = '2026-01-30'
$ws.Cells.Item(57, 1).ClearFormats()
$ws.Cells.Item(57, 2).Value = '17:51:49'
$ws.Cells.Item(57, 3).Value = '17:00'
$ws.Cells.Item(57, 4).Value = 'Bathroom Door'
$ws.Cells.Item(57, 5).Value = 'ENTER'
$ws.Cells.Item(57, 6).Value = 'User ENTERED Bathroom'
# Row 58
$ws.Cells.Item(58, 1).NumberFormat = '@'
$ws.Cells.Item(58, 1).Value = '2026-01-30'
$ws.Cells.Item(58, 1).ClearFormats()
$ws.Cells.Item(58, 2).Value = '17:52:23'
$ws.Cells.Item(58, 3).Value = '17:00'
$ws.Cells.Item(58, 4).Value = 'Bathroom Door'
$ws.Cells.Item(58, 5).Value = 'EXIT'
$ws.Cells.Item(58, 6).Value = 'User EXITED Bathroom'
# Row 59
$ws.Cells.Item(59, 1).NumberFormat = '@'
$ws.Cells.Item(59, 1).Value = '2026-01-30'
$ws.Cells.Item(59, 1).ClearFormats()
$ws.Cells.Item(59, 2).Value = '17:52:40'
$ws.Cells.Item(59, 3).Value = '17:00'
$ws.Cells.Item(59, 4).Value = 'Bathroom Door'
$ws.Cells.Item(59, 5).Value = 'ENTER'
$ws.Cells.Item(59, 6).Value = 'User ENTERED Bathroom'
